# durham_daily_vaccinations.xlsx — "Add files via upload"
#
# The sheet tracks one row per day (columns A-I). The previously "latest"
# row (107, highlighted with the green "Neutral"/shaded style + the
# "daily rate to achieve June 20 target" label in column I) receives a
# corrected dose count, a new day's row is appended after it (taking over
# the highlight + label), and a second blank spacer row is inserted above
# the trailing SUM / difference rows further down the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Make room -----------------------------------------------------
# Insert the new "latest day" row right after the old one (row 108);
# Excel seeds the new row's formatting from the row above (the shaded
# "latest" style), which is exactly what the new row needs.
$ws.Rows("108:108").Insert()

# Insert a second blank spacer row further down (pushes the trailing
# SUM/diff rows from 110/111 to 112/113, and auto-adjusts the formulas
# that reference them).
$ws.Rows("111:111").Insert()

# --- 2. Correct the (now second-to-last) existing data row 107 --------
$ws.Range("B107").Value = 3303

# --- 3. Fill in the new last data row (108) ----------------------------
$ws.Range("A108").Value = 44293
$ws.Range("B108").Value = 3096
$ws.Range("C108").Formula = "=(AVERAGE(B102:B108))"
$ws.Range("D108").Formula = "=(D107-B108)"
$ws.Range("E108").Formula = "=E107+B108"
$ws.Range("F108").Formula = "=D108/C108"
$ws.Range("G108").Formula = "=A108+F108"
$ws.Range("H108").Formula = "=D108/84"
$ws.Range("I108").Value = "daily rate to achieve June 20 target"

# --- 4. Re-color row 107 back to the normal (non-highlighted) style ---
# Row 108 inherited row 107's old "latest row" shading on insert; row 107
# itself needs to drop back to the plain style used by the rest of the
# historical rows (copy formats from row 106, which already has it) and
# lose its "daily rate..." label (that now belongs to row 108).
$ws.Range("A106:H106").Copy()
$ws.Range("A107:H107").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("I107").Clear()

# --- 5. Restore view state ---------------------------------------------
$ws.Range("J115").Select()
$excel.ActiveWindow.ScrollRow = 77

Write-Output "done"
